$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-13 10:21:51", 0.0004),
    @("2023-12-13 10:22:28", 0.0018),
    @("2023-12-13 10:23:31", 0.004200000000000001),
    @("2023-12-13 10:23:54", 0.0008)
)

$startRow = 266
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
